$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.726.92'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.624.56'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.65'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0611'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.34'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0856'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.852.18'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.623.73'
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '65.22'
$ws.Range('E15').Value = '  +1.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.513'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.746.20'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.75'
$ws.Range('E18').Value = '  +8.78%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.76'
$ws.Range('E19').Value = '  +4.55%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0728'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.40'
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.22'
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.16'
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.92'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.05'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  +1.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.451.54'
$ws.Range('E33').Value = '  +7.65%  '
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.568'
$ws.Range('E37').Value = '  -3.25%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0168'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('E39').Value = '  +2.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.97'
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.21'
$ws.Range('E42').Value = '  +2.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.949'
$ws.Range('E43').Value = '  -5.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.764.06'
$ws.Range('E44').Value = '  +2.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.765'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.31'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.49'
$ws.Range('E47').Value = '  +2.74%  '
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0505'
$ws.Range('E49').Value = '  +0.60%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₇0985'
$ws.Range('E50').Value = '  -4.46%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0965'
$ws.Range('E51').Value = '  -1.40%  '
